$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.013.39"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.532.84"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "601.51"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "182.74"
$ws.Range("E6").Value = "  +5.21%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "0.139"
$ws.Range("E9").Value = "  +5.01%  "
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "0.442"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "4.150.78"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "32.38"
$ws.Range("E13").Value = "  +11.69%  "
$ws.Range("D14").Value = "0.136"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "68.016.11"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "3.549.53"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "14.68"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").Value = "400.30"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "8.06"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "74.27"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "5.69"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "6.37"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "7.64"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "24.03"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "1.65"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "163.65"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "4.76"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.923.50"
$ws.Range("E42").Value = "  +3.60%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "6.88"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "27.00"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0744"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "26.99"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "42.51"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "352.51"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0307"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "0.893"
$ws.Range("E51").Value = "  +5.13%  "
